# Auto-generated Excel COM-interop script
# Applies numeric cell updates (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the upstream commit's unified OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 105.22222
$ws.Range("I5").Value = 62.5
$ws.Range("K5").Value = 62.5
$ws.Range("M5").Value = 52.5
$ws.Range("H18").Value = 354
$ws.Range("I18").Value = 354
$ws.Range("K18").Value = 354
$ws.Range("M18").Value = -70
$ws.Range("H21").Value = 20904.715
$ws.Range("I21").Value = 18666.6
$ws.Range("J21").Value = 26500
$ws.Range("K21").Value = 18666.6
$ws.Range("L21").Value = 26500
$ws.Range("M21").Value = -18198.6
$ws.Range("N21").Value = -27436
$ws.Range("H23").Value = 20904.715
$ws.Range("I23").Value = 18666.6
$ws.Range("J23").Value = 26500
$ws.Range("K23").Value = 18666.6
$ws.Range("L23").Value = 26500
$ws.Range("M23").Value = -18432.6
$ws.Range("N23").Value = -26968
$ws.Range("H40").Value = 2269.2307
$ws.Range("I40").Value = 2242.8572
$ws.Range("J40").Value = 2300
$ws.Range("K40").Value = 2242.8572
$ws.Range("L40").Value = 2300
$ws.Range("M40").Value = -2067.8572
$ws.Range("N40").Value = -2650
$ws.Range("H107").Value = 741360.25
$ws.Range("I107").Value = 855238.3
$ws.Range("J107").Value = 1153
$ws.Range("K107").Value = 855238.3
$ws.Range("L107").Value = 1153
$ws.Range("M107").Value = -853318.3
$ws.Range("N107").Value = -4993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3419.9534
$ws.Range("I32").Value = 2370.4
$ws.Range("K32").Value = 2370.4
$ws.Range("M32").Value = -2083.4
$ws.Range("H45").Value = 1142.2858
$ws.Range("I45").Value = 1126.5454
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1126.5454
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -749.5454
$ws.Range("N45").Value = -1954
$ws.Range("H102").Value = 2534.5557
$ws.Range("I102").Value = 2300
$ws.Range("J102").Value = 2651.8333
$ws.Range("K102").Value = 2300
$ws.Range("L102").Value = 2651.8333
$ws.Range("M102").Value = -678
$ws.Range("N102").Value = -5895.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 586.4545000000001
$ws.Range("I107").Value = 584.9091
$ws.Range("J107").Value = 588
$ws.Range("K107").Value = 584.9091
$ws.Range("L107").Value = 588
$ws.Range("M107").Value = 1335.0909
$ws.Range("N107").Value = -4428
$ws.Range("H134").Value = 4984.263
$ws.Range("I134").Value = 2627
$ws.Range("J134").Value = 7105.8
$ws.Range("K134").Value = 7881
$ws.Range("L134").Value = 21317.4
$ws.Range("M134").Value = -5346
$ws.Range("N134").Value = -26387.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1429.08
$ws.Range("I31").Value = 1078.9546
$ws.Range("J31").Value = 3996.6667
$ws.Range("K31").Value = 1078.9546
$ws.Range("L31").Value = 3996.6667
$ws.Range("M31").Value = -783.9546
$ws.Range("N31").Value = -4586.6667
$ws.Range("H34").Value = 1429.08
$ws.Range("I34").Value = 1078.9546
$ws.Range("J34").Value = 3996.6667
$ws.Range("K34").Value = 1078.9546
$ws.Range("L34").Value = 3996.6667
$ws.Range("M34").Value = -876.9546
$ws.Range("N34").Value = -4400.6667
$ws.Range("H51").Value = 21475
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 25300
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 25300
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -26772
$ws.Range("H58").Value = 2393.1538
$ws.Range("I58").Value = 1422.625
$ws.Range("J58").Value = 3946
$ws.Range("K58").Value = 1422.625
$ws.Range("L58").Value = 3946
$ws.Range("M58").Value = -1219.625
$ws.Range("N58").Value = -4352
$ws.Range("H60").Value = 10820.909
$ws.Range("I60").Value = 8000
$ws.Range("K60").Value = 8000
$ws.Range("M60").Value = -7489
$ws.Range("H61").Value = 21475
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 25300
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 25300
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -25996
$ws.Range("H136").Value = 2393.1538
$ws.Range("I136").Value = 1422.625
$ws.Range("J136").Value = 3946
$ws.Range("K136").Value = 4267.875
$ws.Range("L136").Value = 11838
$ws.Range("M136").Value = -1717.875
$ws.Range("N136").Value = -16938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1063.6471
$ws.Range("I5").Value = 633.8182
$ws.Range("J5").Value = 1851.6666
$ws.Range("K5").Value = 1901.4546
$ws.Range("L5").Value = 5554.9998
$ws.Range("M5").Value = -1789.4546
$ws.Range("N5").Value = -5778.9998
$ws.Range("H113").Value = 13889697
$ws.Range("J113").Value = 15625837
$ws.Range("L113").Value = 46877511
$ws.Range("N113").Value = -46881851
$ws.Range("H131").Value = 2994.5908
$ws.Range("J131").Value = 3120.8413
$ws.Range("L131").Value = 9362.5239
$ws.Range("N131").Value = -19442.5239
$ws.Range("H132").Value = 1400
$ws.Range("J132").Value = 1400
$ws.Range("L132").Value = 12600
$ws.Range("N132").Value = -17660
$ws.Range("H135").Value = 1063.6471
$ws.Range("I135").Value = 633.8182
$ws.Range("J135").Value = 1851.6666
$ws.Range("K135").Value = 5704.3638
$ws.Range("L135").Value = 16664.9994
$ws.Range("M135").Value = -3169.3638
$ws.Range("N135").Value = -21734.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 315.3
$ws.Range("J107").Value = 318
$ws.Range("L107").Value = 318
$ws.Range("N107").Value = -4158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1343.4286
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 1467.3334
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 1467.3334
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -1843.3334
$ws.Range("H68").Value = 1300.5
$ws.Range("I68").Value = 1300.5
$ws.Range("K68").Value = 1300.5
$ws.Range("M68").Value = -551.5
$ws.Range("H71").Value = 1300.5
$ws.Range("I71").Value = 1300.5
$ws.Range("K71").Value = 6502.5
$ws.Range("M71").Value = -2758.5
$ws.Range("H136").Value = 6563.278
$ws.Range("I136").Value = 2680.4443
$ws.Range("J136").Value = 10446.111
$ws.Range("K136").Value = 8041.3329
$ws.Range("L136").Value = 31338.333
$ws.Range("M136").Value = -5491.3329
$ws.Range("N136").Value = -36438.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16718002
$ws.Range("I136").Value = 22290110
$ws.Range("J136").Value = 1680
$ws.Range("K136").Value = 66870330
$ws.Range("L136").Value = 5040
$ws.Range("M136").Value = -66867780
$ws.Range("N136").Value = -10140

Write-Host "Applied 173 cell updates across 8 sheets."
